$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the time-range values in column C:
#   8:15-8:20 -> 8:25-8:30
#   8:20-8:25 -> 8:30-8:35
$ws.Range("C2").Value = "8:25-8:30"
$ws.Range("C3").Value = "8:30-8:35"

# Move the active cell selection to C10 (was C8)
$ws.Range("C10").Select()
